$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'66.846.01"
$ws.Cells.Item(2, 5).Value = "  +1.01%  "
$ws.Cells.Item(3, 4).Value = "'3.271.84"
$ws.Cells.Item(3, 5).Value = "  -1.81%  "
$ws.Cells.Item(4, 4).Value = "'0.998"
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$ws.Cells.Item(5, 4).Value = "'576.28"
$ws.Cells.Item(5, 5).Value = "  -1.23%  "
$ws.Cells.Item(6, 4).Value = "'172.28"
$ws.Cells.Item(6, 5).Value = "  -6.86%  "
$ws.Cells.Item(8, 4).Value = "'0.579"
$ws.Cells.Item(8, 5).Value = "  +0.53%  "
$ws.Cells.Item(9, 4).Value = "'3.266.19"
$ws.Cells.Item(9, 5).Value = "  -1.84%  "
$ws.Cells.Item(10, 4).Value = "'0.172"
$ws.Cells.Item(10, 5).Value = "  -4.59%  "
$ws.Cells.Item(11, 4).Value = "'0.570"
$ws.Cells.Item(11, 5).Value = "  -1.74%  "
$ws.Cells.Item(12, 4).Value = "'44.93"
$ws.Cells.Item(12, 5).Value = "  -4.50%  "
$ws.Cells.Item(13, 4).Value = "'0.0000265"
$ws.Cells.Item(13, 5).Value = "  -1.09%  "
$ws.Cells.Item(14, 4).Value = "'686.83"
$ws.Cells.Item(14, 5).Value = "  +1.56%  "
$ws.Cells.Item(15, 4).Value = "'3.803.48"
$ws.Cells.Item(15, 5).Value = "  -1.56%  "
$ws.Cells.Item(16, 4).Value = "'8.24"
$ws.Cells.Item(16, 5).Value = "  -2.69%  "
$ws.Cells.Item(17, 4).Value = "'66.922.48"
$ws.Cells.Item(17, 5).Value = "  +0.93%  "
$ws.Cells.Item(18, 5).Value = "  +1.02%  "
$ws.Cells.Item(19, 4).Value = "'3.271.02"
$ws.Cells.Item(19, 5).Value = "  -1.90%  "
$ws.Cells.Item(20, 4).Value = "'17.19"
$ws.Cells.Item(20, 5).Value = "  -3.83%  "
$ws.Cells.Item(21, 4).Value = "'10.66"
$ws.Cells.Item(21, 5).Value = "  -4.04%  "
$ws.Cells.Item(22, 4).Value = "'0.881"
$ws.Cells.Item(22, 5).Value = "  -1.80%  "
$ws.Cells.Item(23, 4).Value = "'16.80"
$ws.Cells.Item(23, 5).Value = "  -5.57%  "
$ws.Cells.Item(24, 5).Value = "  +2.89%  "
$ws.Cells.Item(25, 4).Value = "'99.43"
$ws.Cells.Item(25, 5).Value = "  -3.05%  "
$ws.Cells.Item(26, 4).Value = "'3.84"
$ws.Cells.Item(26, 5).Value = "  -3.38%  "
$ws.Cells.Item(27, 4).Value = "'2.66"
$ws.Cells.Item(27, 5).Value = "  -3.93%  "
$ws.Cells.Item(28, 4).Value = "'33.43"
$ws.Cells.Item(28, 5).Value = "  +3.24%  "
$ws.Cells.Item(29, 4).Value = "'9.11"
$ws.Cells.Item(29, 5).Value = "  -3.64%  "
$ws.Cells.Item(30, 4).Value = "'8.30"
$ws.Cells.Item(30, 5).Value = "  -2.26%  "
$ws.Cells.Item(31, 5).Value = "  -1.30%  "
$ws.Cells.Item(32, 4).Value = "'568.48"
$ws.Cells.Item(32, 5).Value = "  -6.48%  "
$ws.Cells.Item(35, 5).Value = "  -2.97%  "
$ws.Cells.Item(36, 5).Value = "  -0.01%  "
$ws.Cells.Item(37, 4).Value = "'3.32"
$ws.Cells.Item(37, 5).Value = "  -15.70%  "
$ws.Cells.Item(38, 4).Value = "'55.05"
$ws.Cells.Item(38, 5).Value = "  -1.70%  "
$ws.Cells.Item(39, 5).Value = "  +0.59%  "
$ws.Cells.Item(40, 5).Value = "  +0.36%  "
$ws.Cells.Item(41, 5).Value = "  -4.19%  "
$ws.Cells.Item(42, 4).Value = "'31.36"
$ws.Cells.Item(42, 5).Value = "  -4.13%  "
$ws.Cells.Item(43, 4).Value = "'0.0₃0662"
$ws.Cells.Item(43, 5).Value = "  -5.53%  "
$ws.Cells.Item(44, 4).Value = "'2.97"
$ws.Cells.Item(44, 5).Value = "  -6.91%  "
$ws.Cells.Item(45, 4).Value = "'0.324"
$ws.Cells.Item(45, 5).Value = "  -3.65%  "
$ws.Cells.Item(46, 4).Value = "'0.0401"
$ws.Cells.Item(46, 5).Value = "  -3.21%  "
$ws.Cells.Item(49, 5).Value = "  -1.18%  "
$ws.Cells.Item(50, 4).Value = "'1.36"
$ws.Cells.Item(50, 5).Value = "  +4.28%  "
$ws.Cells.Item(51, 4).Value = "'129.51"
$ws.Cells.Item(51, 5).Value = "  -1.23%  "

# Row 33/34 swap: Maker <-> Cosmos (with updated prices)
$ws.Cells.Item(33, 2).Value = "Cosmos"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(33, 4).Value = "'10.78"
$ws.Cells.Item(33, 5).Value = "  -2.37%  "
$ws.Cells.Item(34, 2).Value = "Maker"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(34, 4).Value = "'3.833.57"
$ws.Cells.Item(34, 5).Value = "  +0.66%  "

# Row 47/48 swap: Stellar <-> FirstDigitalUSD (with updated prices)
$ws.Cells.Item(47, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(47, 4).Value = "'1.00"
$ws.Cells.Item(47, 5).Value = "  -0.13%  "
$ws.Cells.Item(48, 2).Value = "Stellar"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(48, 4).Value = "'0.126"
$ws.Cells.Item(48, 5).Value = "  -1.04%  "
